$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16
$ws.Range("A16").Value = "IRST015"
$ws.Range("B16").Value = "Right"
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 2

# Row 17
$ws.Range("A17").Value = "IRST016"
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 2.2
$ws.Range("F17").Value = 2.2

$ws.Range("H16").Value = "Upper outer; unknown size"
$ws.Range("B17").Value = "Right "

$ws.Range("B17").Select()
